$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A34").Value = "محمد تحسين طاهر"
$ws.Range("B34").Value = "مدري وش يعني معامله"

$null = $ws.Range("B34").Select()
